# Auto-generated edit script: refreshes the 3 schedule sheets
# (LP1912, LP1912-215, 6203-6173) with the 07:44:08 scrape data,
# matching the commit 'Horarios actualizados Linea 141 - 427'.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Última actualización: 07:44:08"
$ws.Range("A3").Value = "Total filas: 98"
$rows = @(
    @("00:04:05","01:12","215_ALUAR",68,"LP1912"),
    @("00:04:05","01:58","14_ABASTO",114,"LP1912"),
    @("01:07:17","02:58","215_ALUAR",111,"LP1912"),
    @("01:49:17","03:03","215_ALUAR",74,"LP1912"),
    @("01:49:17","03:48","14_ABASTO",119,"LP1912"),
    @("02:42:39","03:54","14_ABASTO",72,"LP1912"),
    @("02:13:28","03:55","14_ABASTO",102,"LP1912"),
    @("02:13:28","04:01","81_EL PELIGRO",108,"LP1912"),
    @("03:23:38","04:45","215A_EL PATO",82,"LP1912"),
    @("02:56:05","04:46","215A_EL PATO",110,"LP1912"),
    @("02:56:05","04:53","11_ETCHEVERRY",117,"LP1912"),
    @("04:53:50","04:54","11_ETCHEVERRY",1,"LP1912"),
    @("04:53:50","05:15","14_ABASTO",22,"LP1912"),
    @("03:23:38","05:16","17_ROMERO",113,"LP1912"),
    @("04:40:32","05:17","14_ABASTO",37,"LP1912"),
    @("05:18:23","05:20","14_ABASTO",2,"LP1912"),
    @("03:23:38","05:22","23_HERNANDEZ",119,"LP1912"),
    @("03:54:40","05:34","215B_EL PATO",100,"LP1912"),
    @("04:18:06","05:35","215B_EL PATO",77,"LP1912"),
    @("04:18:06","05:36","14_ABASTO",78,"LP1912"),
    @("03:54:40","05:39","14_ABASTO",105,"LP1912"),
    @("03:54:40","05:46","15_ABASTO",112,"LP1912"),
    @("05:47:32","05:47","15_ABASTO",0,"LP1912"),
    @("05:47:32","05:49","14_ABASTO",2,"LP1912"),
    @("04:40:32","06:04","16_SANTA ANA",84,"LP1912"),
    @("06:02:16","06:05","16_SANTA ANA",3,"LP1912"),
    @("04:18:06","06:09","16_SANTA ANA",111,"LP1912"),
    @("04:40:32","06:11","215A_EL PATO",91,"LP1912"),
    @("04:18:06","06:12","215A_EL PATO",114,"LP1912"),
    @("04:18:06","06:14","225_HARAS DEL SUR",116,"LP1912"),
    @("04:40:32","06:21","26_HERNANDEZ",101,"LP1912"),
    @("04:40:32","06:27","23_HERNANDEZ",107,"LP1912"),
    @("04:40:32","06:29","86_EST CHICA-ESC AGRARIA",109,"LP1912"),
    @("06:02:16","06:30","86_EST CHICA-ESC AGRARIA",28,"LP1912"),
    @("04:40:32","06:31","16_SANTA ANA",111,"LP1912"),
    @("04:53:50","06:44","225_C ROCA-H SUR",111,"LP1912"),
    @("04:53:50","06:46","215C_EL PATO",113,"LP1912"),
    @("05:18:23","06:58","10_OLMOS",100,"LP1912"),
    @("05:18:23","06:59","14_ABASTO",101,"LP1912"),
    @("06:02:16","07:00","14_ABASTO",58,"LP1912"),
    @("06:37:24","07:01","16_SANTA ANA",24,"LP1912"),
    @("05:47:32","07:04","23_HERNANDEZ",77,"LP1912"),
    @("05:18:23","07:05","15_ABASTO",107,"LP1912"),
    @("06:02:16","07:05","23_HERNANDEZ",63,"LP1912"),
    @("05:18:23","07:07","225_GOMEZ",109,"LP1912"),
    @("05:18:23","07:11","215A_EL PATO",113,"LP1912"),
    @("07:14:27","07:14","11_ETCHEVERRY",0,"LP1912"),
    @("05:18:23","07:15","11_ETCHEVERRY",117,"LP1912"),
    @("06:37:24","07:16","16_SANTA ANA",39,"LP1912"),
    @("06:02:16","07:16","11_ETCHEVERRY",74,"LP1912"),
    @("05:47:32","07:21","26_HERNANDEZ",94,"LP1912"),
    @("06:02:16","07:23","10_OLMOS",81,"LP1912"),
    @("05:47:32","07:27","10_OLMOS",100,"LP1912"),
    @("05:47:32","07:31","11_ETCHEVERRY",104,"LP1912"),
    @("05:47:32","07:31","16_SANTA ANA",104,"LP1912"),
    @("06:02:16","07:32","11_ETCHEVERRY",90,"LP1912"),
    @("05:47:32","07:32","84_COLONIA URQUIZA-ESC 49",105,"LP1912"),
    @("06:37:24","07:34","23_HERNANDEZ",57,"LP1912"),
    @("05:47:32","07:36","27_EL RETIRO",109,"LP1912"),
    @("07:14:27","07:37","23_HERNANDEZ",23,"LP1912"),
    @("06:02:16","07:37","27_EL RETIRO",95,"LP1912"),
    @("05:47:32","07:39","10_OLMOS",112,"LP1912"),
    @("07:44:08","07:45","14_ABASTO",1,"LP1912"),
    @("06:37:24","07:47","14_ABASTO",70,"LP1912"),
    @("06:02:16","07:48","14_ABASTO",106,"LP1912"),
    @("06:02:16","07:51","215D_EL PATO",109,"LP1912"),
    @("07:44:08","07:52","215D_EL PATO",8,"LP1912"),
    @("07:44:08","07:55","10_OLMOS",11,"LP1912"),
    @("07:14:27","07:58","16_SANTA ANA",44,"LP1912"),
    @("07:14:27","08:03","11_ETCHEVERRY",49,"LP1912"),
    @("07:44:08","08:04","11_ETCHEVERRY",20,"LP1912"),
    @("07:44:08","08:11","16_SANTA ANA",27,"LP1912"),
    @("06:37:24","08:12","15_ABASTO",95,"LP1912"),
    @("07:44:08","08:13","10_OLMOS",29,"LP1912"),
    @("06:37:24","08:21","26_HERNANDEZ",104,"LP1912"),
    @("06:37:24","08:22","16_P MOR-SANTA ANA",105,"LP1912"),
    @("06:37:24","08:23","215B_EL PATO",106,"LP1912"),
    @("07:44:08","08:23","16_P MOR-SANTA ANA",39,"LP1912"),
    @("06:37:24","08:27","84_COLONIA URQUIZA-ESC 49",110,"LP1912"),
    @("06:52:38","08:42","81_EL PELIGRO",110,"LP1912"),
    @("07:14:27","08:43","14_ABASTO",89,"LP1912"),
    @("07:44:08","08:44","14_ABASTO",60,"LP1912"),
    @("07:14:27","08:54","17_ROMERO",100,"LP1912"),
    @("07:14:27","09:01","215A_EL PATO",107,"LP1912"),
    @("07:44:08","09:02","215A_EL PATO",78,"LP1912"),
    @("07:14:27","09:07","23_HERNANDEZ",113,"LP1912"),
    @("07:44:08","09:08","23_HERNANDEZ",84,"LP1912"),
    @("07:14:27","09:10","16_P MOR-SANTA ANA",116,"LP1912"),
    @("07:44:08","09:11","16_P MOR-SANTA ANA",87,"LP1912"),
    @("07:44:08","09:14","16_SANTA ANA",90,"LP1912"),
    @("07:44:08","09:17","27_EL RETIRO",93,"LP1912"),
    @("07:44:08","09:21","26_HERNANDEZ",97,"LP1912"),
    @("07:44:08","09:23","17_ROMERO",99,"LP1912"),
    @("07:44:08","09:24","11_ETCHEVERRY",100,"LP1912"),
    @("07:44:08","09:32","15_ABASTO",108,"LP1912"),
    @("07:44:08","09:33","10_OLMOS",109,"LP1912"),
    @("07:44:08","09:36","23_HERNANDEZ",112,"LP1912"),
    @("07:44:08","09:42","215C_EL PATO",118,"LP1912")
)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 6
    $ws.Cells.Item($rowNum, 1).Value = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}

# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Última actualización: 07:44:08"
$ws.Range("A3").Value = "Total filas: 17"
$rows = @(
    @("00:04:05","01:12","215_ALUAR",68,"LP1912"),
    @("01:07:17","02:58","215_ALUAR",111,"LP1912"),
    @("01:49:17","03:03","215_ALUAR",74,"LP1912"),
    @("03:23:38","04:45","215A_EL PATO",82,"LP1912"),
    @("02:56:05","04:46","215A_EL PATO",110,"LP1912"),
    @("03:54:40","05:34","215B_EL PATO",100,"LP1912"),
    @("04:18:06","05:35","215B_EL PATO",77,"LP1912"),
    @("04:40:32","06:11","215A_EL PATO",91,"LP1912"),
    @("04:18:06","06:12","215A_EL PATO",114,"LP1912"),
    @("04:53:50","06:46","215C_EL PATO",113,"LP1912"),
    @("05:18:23","07:11","215A_EL PATO",113,"LP1912"),
    @("06:02:16","07:51","215D_EL PATO",109,"LP1912"),
    @("07:44:08","07:52","215D_EL PATO",8,"LP1912"),
    @("06:37:24","08:23","215B_EL PATO",106,"LP1912"),
    @("07:14:27","09:01","215A_EL PATO",107,"LP1912"),
    @("07:44:08","09:02","215A_EL PATO",78,"LP1912"),
    @("07:44:08","09:42","215C_EL PATO",118,"LP1912")
)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 6
    $ws.Cells.Item($rowNum, 1).Value = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}

# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Última actualización: 07:44:08"
$ws.Range("A3").Value = "Total filas: 17"
$rows = @(
    @("00:04:05","00:05","215A_LA PLATA",1,"L6173"),
    @("03:54:40","05:43","215A_LA PLATA",109,"L6173"),
    @("04:18:06","05:44","215A_LA PLATA",86,"L6173"),
    @("04:53:50","06:08","215A_LA PLATA",75,"L6173"),
    @("04:18:06","06:09","215A_LA PLATA",111,"L6173"),
    @("04:53:50","06:32","215C_LA PLATA",99,"L6203"),
    @("04:40:32","06:33","215C_LA PLATA",113,"L6203"),
    @("05:18:23","07:00","215B_LP-P MOR-1 Y 57",102,"L6173"),
    @("05:47:32","07:35","215A_LA PLATA",108,"L6173"),
    @("06:52:38","08:06","215C_LA PLATA",74,"L6203"),
    @("06:37:24","08:07","215C_LA PLATA",90,"L6203"),
    @("07:14:27","08:14","215C_LA PLATA",60,"L6203"),
    @("07:44:08","08:17","215C_LA PLATA",33,"L6203"),
    @("06:37:24","08:30","215A_LA PLATA",113,"L6173"),
    @("06:52:38","08:34","215A_LA PLATA",102,"L6173"),
    @("07:14:27","08:35","215A_LA PLATA",81,"L6173"),
    @("07:14:27","09:09","215D_LA PLATA",115,"L6203")
)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 6
    $ws.Cells.Item($rowNum, 1).Value = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}
